# Slide 20 ("Compound Statement"), Content Placeholder 2 (shape 2): split the
# first paragraph's runs so that:
#   "A sequence of zero or more statements enclosed in "  -> run 1 (unchanged)
#   "braces \u201c"                                        -> run 2 (new split)
#   "{"  (Consolas)                                        -> run 3 (unchanged, Consolas)
#   "\u201d "                                               -> run 4 (new split)
#   "and \u201c"                                            -> run 5 (new split)
#   "}"  (Consolas)                                        -> run 6 (unchanged, Consolas)
#   "\u201d."                                               -> run 7 (unchanged)
# The visible text is left exactly as it was; only the run boundaries change,
# matching the target edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Split off "A sequence of zero or more statements enclosed in " (chars 1-50)
# from "braces “" (chars 51-58).
$tr.Characters(1, 50).Text = "A sequence of zero or more statements enclosed in "

# Rewrite "braces “" (chars 51-58) as its own run, separate from the "{" run
# that follows it (char 59).
$tr.Characters(51, 8).Text = "braces “"

# Split "” and “" (chars 60-66) into "” " (60-61) and "and “" (62-66).
$tr.Characters(60, 2).Text = "” "
$tr.Characters(62, 5).Text = "and “"
